# LMS-2391 Split TSV files where necessary.
# Duplicate the OD600 "MGP1" row into a new row 6 on the openbis-data sheet,
# switch the active tab from openbis-metadata to openbis-data, update the
# selection on openbis-data to the newly-added row, and set the sheet's
# page setup (paper size / orientation) to match.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("openbis-data")

# Duplicate row 2 ("MGP1"/"OD600" values) into row 6, preserving shared
# string reuse and exact numeric values by copying cell-by-cell.
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $src = $wsData.Cells.Item(2, $col)
    $dst = $wsData.Cells.Item(6, $col)
    $dst.Value = $src.Text
}

# Page setup for the data sheet.
$wsData.PageSetup.PaperSize = 10
$wsData.PageSetup.Orientation = 1

# Move the active tab to openbis-data (this also moves tabSelected from
# openbis-metadata's sheetView to openbis-data's sheetView).
$wsData.Activate()

# Update the selection on openbis-data to the newly added row.
[void]$wsData.Range("A6:XFD6").Select()
